# ITS-incident-reporting.xlsx — "ADD RTS and ITS incident reporting"
#
# Content changes:
#  - library_content sheet: rename library/framework name, swap provider/copyright
#    from "ESMA" to "ESA"
#  - ITS-incident-reporting sheet: renumber sub-items "1.a/1.b/1.c" -> "1.1.a/1.1.b/1.1.c"
#  - Switch the active sheet/selection from ITS-incident-reporting to library_content

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("library_content")
$ws2 = $wb.Worksheets.Item("ITS-incident-reporting")

# --- library_content updates -------------------------------------------------
# (Ordered to match the shared-string append order of the source edit: the new
# "ITS DORA incident reporting" label comes first, then the renumbered "1.1.x"
# sub-items, then the "ESA" relabeling.)
$ws1.Range("B5").Value2 = "ITS DORA incident reporting"     # library_name
$ws1.Range("B12").Value2 = "ITS DORA incident reporting"    # framework_name

# --- ITS-incident-reporting updates ------------------------------------------
$ws2.Range("C14").Value2 = "1.1.a"
$ws2.Range("C15").Value2 = "1.1.b"
$ws2.Range("C16").Value2 = "1.1.c"

# --- library_content updates (continued) -------------------------------------
$ws1.Range("B7").Value2 = "ESA"                              # library_copyright
$ws1.Range("B8").Value2 = "ESA"                              # library_provider

# --- view / selection state ---------------------------------------------------
# Move the scroll position and selection on ITS-incident-reporting, then make
# library_content the active/visible tab with its own selection.
$ws2.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 11
$aw.ScrollColumn = 1
$ws2.Range("G15").Select()

$ws1.Activate()
$ws1.Range("C6").Select()
